$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "80-54="
$t.Cell(1,2).Range.Text = "77-70="
$t.Cell(1,3).Range.Text = "63+5="
$t.Cell(1,4).Range.Text = "48+42="
$t.Cell(1,5).Range.Text = "75-22="
$t.Cell(2,1).Range.Text = "45+29="
$t.Cell(2,2).Range.Text = "31+37="
$t.Cell(2,3).Range.Text = "49-28="
$t.Cell(2,4).Range.Text = "33+53="
$t.Cell(2,5).Range.Text = "86-70="
$t.Cell(3,1).Range.Text = "16-0="
$t.Cell(3,2).Range.Text = "16+20="
$t.Cell(3,3).Range.Text = "14-0="
$t.Cell(3,4).Range.Text = "3+13="
$t.Cell(3,5).Range.Text = "53-13="
$t.Cell(4,1).Range.Text = "10+70="
$t.Cell(4,2).Range.Text = "51-1="
$t.Cell(4,3).Range.Text = "77-36="
$t.Cell(4,4).Range.Text = "2+9="
$t.Cell(4,5).Range.Text = "10+18="
$t.Cell(5,1).Range.Text = "38-34="
$t.Cell(5,2).Range.Text = "84-23="
$t.Cell(5,3).Range.Text = "39-2="
$t.Cell(5,4).Range.Text = "36-1="
$t.Cell(5,5).Range.Text = "86-80="
$t.Cell(6,1).Range.Text = "85-70="
$t.Cell(6,2).Range.Text = "44+25="
$t.Cell(6,3).Range.Text = "12+31="
$t.Cell(6,4).Range.Text = "69-15="
$t.Cell(6,5).Range.Text = "14+21="
$t.Cell(7,1).Range.Text = "28+35="
$t.Cell(7,2).Range.Text = "78-44="
$t.Cell(7,3).Range.Text = "60+3="
$t.Cell(7,4).Range.Text = "14+70="
$t.Cell(7,5).Range.Text = "45-7="
$t.Cell(8,1).Range.Text = "57+23="
$t.Cell(8,2).Range.Text = "74-33="
$t.Cell(8,3).Range.Text = "22+72="
$t.Cell(8,4).Range.Text = "56-10="
$t.Cell(8,5).Range.Text = "46-30="
$t.Cell(9,1).Range.Text = "11+45="
$t.Cell(9,2).Range.Text = "27-8="
$t.Cell(9,3).Range.Text = "53+7="
$t.Cell(9,4).Range.Text = "82-63="
$t.Cell(9,5).Range.Text = "37+34="
$t.Cell(10,1).Range.Text = "0+38="
$t.Cell(10,2).Range.Text = "30+30="
$t.Cell(10,3).Range.Text = "23+75="
$t.Cell(10,4).Range.Text = "22-16="
$t.Cell(10,5).Range.Text = "51-13="
$t.Cell(11,1).Range.Text = "15+42="
$t.Cell(11,2).Range.Text = "76-9="
$t.Cell(11,3).Range.Text = "84+13="
$t.Cell(11,4).Range.Text = "97-83="
$t.Cell(11,5).Range.Text = "33+60="
$t.Cell(12,1).Range.Text = "45+45="
$t.Cell(12,2).Range.Text = "71+26="
$t.Cell(12,3).Range.Text = "58-21="
$t.Cell(12,4).Range.Text = "37+29="
$t.Cell(12,5).Range.Text = "40-24="
$t.Cell(13,1).Range.Text = "85-67="
$t.Cell(13,2).Range.Text = "21-4="
$t.Cell(13,3).Range.Text = "34+40="
$t.Cell(13,4).Range.Text = "46+2="
$t.Cell(13,5).Range.Text = "22+62="
$t.Cell(14,1).Range.Text = "14+53="
$t.Cell(14,2).Range.Text = "47+34="
$t.Cell(14,3).Range.Text = "40+17="
$t.Cell(14,4).Range.Text = "51-27="
$t.Cell(14,5).Range.Text = "99-87="
$t.Cell(15,1).Range.Text = "51-44="
$t.Cell(15,2).Range.Text = "41+14="
$t.Cell(15,3).Range.Text = "55-4="
$t.Cell(15,4).Range.Text = "37-23="
$t.Cell(15,5).Range.Text = "35+61="
$t.Cell(16,1).Range.Text = "76-62="
$t.Cell(16,2).Range.Text = "65+10="
$t.Cell(16,3).Range.Text = "93-40="
$t.Cell(16,4).Range.Text = "35-25="
$t.Cell(16,5).Range.Text = "57+2="
$t.Cell(17,1).Range.Text = "29+10="
$t.Cell(17,2).Range.Text = "66-62="
$t.Cell(17,3).Range.Text = "25-1="
$t.Cell(17,4).Range.Text = "51-1="
$t.Cell(17,5).Range.Text = "48-43="
$t.Cell(18,1).Range.Text = "98-83="
$t.Cell(18,2).Range.Text = "70-34="
$t.Cell(18,3).Range.Text = "42-16="
$t.Cell(18,4).Range.Text = "42-15="
$t.Cell(18,5).Range.Text = "2+48="
$t.Cell(19,1).Range.Text = "0+6="
$t.Cell(19,2).Range.Text = "74+17="
$t.Cell(19,3).Range.Text = "28+37="
$t.Cell(19,4).Range.Text = "32+57="
$t.Cell(19,5).Range.Text = "36-13="
$t.Cell(20,1).Range.Text = "81-46="
$t.Cell(20,2).Range.Text = "92-17="
$t.Cell(20,3).Range.Text = "16+64="
$t.Cell(20,4).Range.Text = "72-34="
$t.Cell(20,5).Range.Text = "91-20="
